# Re-applies the corrected fixture-to-row assignment for rows 188-190 (3-way
# rotation) and rows 192-193 (swap) in the "Austria Bundesliga" sheet, as
# captured by the 08-05-2024 20:15 base update.
#
# Each affected row keeps its position (A id, C Div, D Date, G/H/I scores,
# V/W/X profit columns all stay put) but the match identity (id in col B),
# the two teams (E/F) and the odds (J:U) move to reflect the corrected
# fixture list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow {
    param(
        [int]$Row,
        [string]$Id,
        [string]$HomeTeam,
        [string]$AwayTeam,
        [double]$OddHOp,
        [double]$OddDOp,
        [double]$OddAOp,
        [double]$OddH,
        [double]$OddD,
        [double]$OddA,
        [double]$Ah,
        [double]$OddAHH,
        [double]$OddAHA,
        [double]$AhOU,
        [object]$OddAHOver,
        [object]$OddAHUnder
    )

    # Column B ("id") is stored as text even though it looks numeric - force
    # text entry the same way Excel does for a quote-prefixed value, then
    # strip the quote-prefix style back off so no visible formatting sticks.
    $ws.Cells.Item($Row, 2).Value = "'" + $Id
    $ws.Cells.Item($Row, 2).Style = "Normal"

    $ws.Cells.Item($Row, 5).Value = $HomeTeam   # E - HomeTeam
    $ws.Cells.Item($Row, 6).Value = $AwayTeam   # F - AwayTeam

    $ws.Cells.Item($Row, 10).Value = $OddHOp    # J - oddH_op
    $ws.Cells.Item($Row, 11).Value = $OddDOp    # K - oddD_op
    $ws.Cells.Item($Row, 12).Value = $OddAOp    # L - oddA_op
    $ws.Cells.Item($Row, 13).Value = $OddH      # M - oddH
    $ws.Cells.Item($Row, 14).Value = $OddD      # N - oddD
    $ws.Cells.Item($Row, 15).Value = $OddA      # O - oddA
    $ws.Cells.Item($Row, 16).Value = $Ah        # P - Ah
    $ws.Cells.Item($Row, 17).Value = $OddAHH    # Q - oddAHH
    $ws.Cells.Item($Row, 18).Value = $OddAHA    # R - oddAHA
    $ws.Cells.Item($Row, 19).Value = $AhOU      # S - AhOU

    if ($null -ne $OddAHOver) {
        $ws.Cells.Item($Row, 20).Value = $OddAHOver   # T - oddAHOver
    }
    if ($null -ne $OddAHUnder) {
        $ws.Cells.Item($Row, 21).Value = $OddAHUnder  # U - oddAHUnder
    }
}

Set-MatchRow 188 "7948287" "Austria Lustenau" "FC Blau Weiss Linz" `
    3.25 3.2 2.2 3.4 3.2 2.15 0.25 1.975 1.875 2.5 2.05 1.8

Set-MatchRow 189 "7948286" "FK Austria Vienna" "Wolfsberger AC" `
    1.727 3.6 4.333 1.8 3.6 4.333 -0.75 2.025 1.825 2.5 1.9 1.95

Set-MatchRow 190 "7948288" "WSG Swarovski Tirol" "SCR Altach" `
    2.4 3.1 2.875 2.75 3.1 2.625 0 1.975 1.875 2.25 $null $null

Set-MatchRow 192 "7947240" "LASK Linz" "SK Sturm Graz" `
    2.6 3.4 2.5 2.7 3.4 2.5 0 2 1.85 2.5 1.975 1.875

Set-MatchRow 193 "7948261" "Austria Klagenfurt" "Rapid Vienna" `
    3.2 3.5 2.1 3.5 3.6 2 0.5 1.85 2 2.75 1.95 1.9
